# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing columns (B..G) and populating H2:H14 with the
# save indicator values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (bold, bordered, centered) onto H1, then
# set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the data rows (H2:H14) with the save values.
$values = @(1, 0, 1, 0, 1, 1, 0, 1, 1, 0, 0, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
